$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AA2").Value = 13
$ws.Range("AC2").Value = 10
$ws.Range("AL2").Value = 29
$ws.Range("AN2").Value = 23

# Row 10
$ws.Range("G10").Value = 1.8
$ws.Range("H10").Value = 3.1
$ws.Range("I10").Value = 5.25
$ws.Range("J10").Value = 2.5
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 2.63
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57
$ws.Range("S10").Value = 3.55
$ws.Range("T10").Value = 1.29
$ws.Range("U10").Value = 4.5
$ws.Range("V10").Value = 1.18
$ws.Range("W10").Value = 1.53
$ws.Range("X10").Value = 2.38
$ws.Range("Y10").Value = 2.2
$ws.Range("Z10").Value = 1.62
$ws.Range("AA10").Value = 5.5
$ws.Range("AB10").Value = 7.5
$ws.Range("AE10").Value = 17
$ws.Range("AG10").Value = 7
$ws.Range("AL10").Value = 11
$ws.Range("AM10").Value = 23
$ws.Range("AN10").Value = 17
$ws.Range("AR10").Value = 1.8
$ws.Range("AS10").Value = 2.05

# Row 11
$ws.Range("G11").Value = 2.1
$ws.Range("H11").Value = 3.05
$ws.Range("I11").Value = 3.75
$ws.Range("J11").Value = 2.8
$ws.Range("K11").Value = 1.9
$ws.Range("L11").Value = 4.4
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 1.47
$ws.Range("P11").Value = 2.57
$ws.Range("Q11").Value = 2.37
$ws.Range("R11").Value = 1.55
$ws.Range("S11").Value = 3.8
$ws.Range("T11").Value = 1.26
$ws.Range("U11").Value = 4.25
$ws.Range("V11").Value = 1.2
$ws.Range("W11").Value = 1.57
$ws.Range("X11").Value = 2.3
$ws.Range("Y11").Value = 2.02
$ws.Range("Z11").Value = 1.7
$ws.Range("AA11").Value = 5.9
$ws.Range("AB11").Value = 9.5
$ws.Range("AC11").Value = 9.5
$ws.Range("AD11").Value = 21
$ws.Range("AE11").Value = 21
$ws.Range("AF11").Value = 40
$ws.Range("AG11").Value = 6
$ws.Range("AH11").Value = 6.3
$ws.Range("AI11").Value = 19
$ws.Range("AJ11").Value = 120
$ws.Range("AL11").Value = 8.5
$ws.Range("AM11").Value = 20
$ws.Range("AN11").Value = 14.5
$ws.Range("AO11").Value = 65
$ws.Range("AP11").Value = 45
$ws.Range("AQ11").Value = 65
$ws.Range("AR11").Value = 1.81
$ws.Range("AS11").Value = 1.97

# Row 12
$ws.Range("G12").Value = 1.14
$ws.Range("H12").Value = 6.9
$ws.Range("J12").Value = 1.5
$ws.Range("K12").Value = 2.75
$ws.Range("L12").Value = 13
$ws.Range("M12").Value = 1.03
$ws.Range("N12").Value = 10.25
$ws.Range("O12").Value = 1.16
$ws.Range("P12").Value = 4.75
$ws.Range("Q12").Value = 1.52
$ws.Range("R12").Value = 2.45
$ws.Range("U12").Value = 2.25
$ws.Range("V12").Value = 1.6
$ws.Range("W12").Value = 1.28
$ws.Range("X12").Value = 3.45
$ws.Range("Y12").Value = 2.55
$ws.Range("Z12").Value = 1.45
$ws.Range("AA12").Value = 6.3
$ws.Range("AB12").Value = 5.7
$ws.Range("AC12").Value = 11.75
$ws.Range("AG12").Value = 10.25
$ws.Range("AH12").Value = 16.5
$ws.Range("AI12").Value = 45

# Row 16
$ws.Range("M16").Value = 1.1
$ws.Range("N16").Value = 7
$ws.Range("AL16").Value = 10

# Row 17
$ws.Range("G17").Value = 2.15
$ws.Range("I17").Value = 2.8
$ws.Range("J17").Value = 2.75
$ws.Range("L17").Value = 3.25
$ws.Range("N17").Value = 17
$ws.Range("AA17").Value = 12
$ws.Range("AC17").Value = 9.5
$ws.Range("AD17").Value = 23
$ws.Range("AE17").Value = 17
$ws.Range("AL17").Value = 13
$ws.Range("AM17").Value = 17
$ws.Range("AN17").Value = 11
$ws.Range("AO17").Value = 29
$ws.Range("AP17").Value = 21
$ws.Range("AQ17").Value = 23

# Row 21
$ws.Range("G21").Value = 2.05
$ws.Range("H21").Value = 3.6
$ws.Range("I21").Value = 3
$ws.Range("J21").Value = 2.63
$ws.Range("K21").Value = 2.4
$ws.Range("L21").Value = 3.5
$ws.Range("M21").Value = 1.02
$ws.Range("N21").Value = 19
$ws.Range("O21").Value = 1.14
$ws.Range("P21").Value = 5.5
$ws.Range("Q21").Value = 1.53
$ws.Range("R21").Value = 2.4
$ws.Range("S21").Value = 1.83
$ws.Range("T21").Value = 2.03
$ws.Range("U21").Value = 2.2
$ws.Range("V21").Value = 1.62
$ws.Range("W21").Value = 1.29
$ws.Range("X21").Value = 3.5
$ws.Range("Y21").Value = 1.44
$ws.Range("Z21").Value = 2.63
$ws.Range("AB21").Value = 13
$ws.Range("AC21").Value = 9
$ws.Range("AD21").Value = 21
$ws.Range("AE21").Value = 15
$ws.Range("AF21").Value = 19
$ws.Range("AG21").Value = 19
$ws.Range("AI21").Value = 11
$ws.Range("AL21").Value = 15
$ws.Range("AM21").Value = 19
$ws.Range("AN21").Value = 12
$ws.Range("AO21").Value = 34
$ws.Range("AP21").Value = 21
$ws.Range("AQ21").Value = 23

# Row 22
$ws.Range("S22").Value = 1.95
$ws.Range("T22").Value = 1.9

# Row 24
$ws.Range("Q24").Value = 2.5
$ws.Range("R24").Value = 1.5
$ws.Range("AR24").Value = 1.93
$ws.Range("AS24").Value = 1.93

# Row 27
$ws.Range("AA27").Value = 7.8
$ws.Range("AB27").Value = 9.25
$ws.Range("AC27").Value = 8.25

# Row 30
$ws.Range("G30").Value = 2.19
$ws.Range("H30").Value = 3.05
$ws.Range("I30").Value = 3.35

Write-Host "Applied 156 cell updates across rows 2, 10, 11, 12, 16, 17, 21, 22, 24, 27, 30"
